# "clean up large accuracy check update"
#
# The libraryDate column (A) had typo'd / inconsistent dates
# (09.6.19, 09.6.20, 09.6.21, 09.6.22) for what is really a single
# library prep date. Normalize every row to the corrected value
# "09.06.19", and leave the selection parked below the data (A28) the
# way the workbook was left after the cleanup pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The column is General-formatted text that merely looks like a date
# ("09.06.19"); assigning that string straight to .Value would make
# Excel's smart-type detection reinterpret it as an actual date
# serial. Briefly force the range to text ("@") while we write the
# corrected values, then restore "General" so the stored number
# format matches the original workbook exactly.
$dateRange = $ws.Range("A2:A27")
$dateRange.NumberFormat = "@"

for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 1).Value = "09.06.19"
}

$dateRange.NumberFormat = "General"

# Selection left at A28 (just past the last data row) after the cleanup.
$ws.Range("A28").Select()
